# Fixed harvester column in rnaSamples -- holly added S.GISH to harvester in bioSamples
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("harvester") currently holds "Retrofitted_2002" for every data row (2-49).
# Holly's correction replaces that value with "S.GISH" for all of those rows.
$ws.Range("B2:B49").Value = "S.GISH"
